# Add 2022-Q3 data as a new sheet, keep 2022-Q2 sheet, and update the
# "总计" (totals) sheet with a new summary row for 2022-Q3.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate the existing "2022-Q2" sheet so its content survives
# at the end of the workbook (it becomes sheetId=3), then reuse the
# original "2022-Q2" sheet object (sheetId=2) as the new "2022-Q3" sheet.
# This matches the target sheetId/order: 总计(1), 2022-Q3(2), 2022-Q2(3).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item(2)
$q2.Copy($null, $q2)
$q2Copy = $wb.Worksheets.Item(3)
$q2Copy.Name = "zz_temp_q2"
$q2.Name = "2022-Q3"
$q2Copy.Name = "2022-Q2"

$q3 = $wb.Worksheets.Item(2)
$total = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# Helper: write a value into a cell while forcing it to be stored as
# TEXT (even when it looks like a number, e.g. "0.66" or "007254") and
# without leaving a NumberFormat-driven style behind on the cell. We do
# this by writing through a scratch cell formatted as Text, then using
# PasteSpecial(values-only) to transplant just the typed value.
# ---------------------------------------------------------------------
function Set-TextValue($sheet, $row, $col, $val) {
    $stage = $sheet.Cells.Item(500, 500)
    $stage.NumberFormat = "@"
    $stage.Value2 = $val
    $stage.Copy()
    $target = $sheet.Cells.Item($row, $col)
    $target.PasteSpecial(-4163)
    $stage.Clear()
}

# ---------------------------------------------------------------------
# Step 2: clear the (now) "2022-Q3" sheet and rebuild its content.
# ---------------------------------------------------------------------
$q3.Cells.Clear()

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    Set-TextValue $q3 1 $col $h
    $col = $col + 1
}

$rows = @(
    @("007254", "广发均衡价值混合", "0.66", "89.74", "5.69", "0.0376", 4),
    @("011003", "同泰大健康主题混合C", "0.28", "94.09", "5.96", "0.0167", 1),
    @("011002", "同泰大健康主题混合A", "0.12", "94.09", "5.96", "0.0072", 1),
    @("008842", "同泰远见灵活配置混合A", "0.18", "93.90", "3.42", "0.0062", 3),
    @("014285", "鑫元健康产业混合A", "0.12", "78.73", "3.13", "0.0038", 9),
    @("008843", "同泰远见灵活配置混合C", "0.07", "93.90", "3.42", "0.0024", 3),
    @("014286", "鑫元健康产业混合C", "0.07", "78.73", "3.13", "0.0022", 9),
    @("006689", "方正富邦信泓灵活配置混合A", "0.03", "93.46", "4.40", "0.0013", 10),
    @("008182", "方正富邦信泓灵活配置混合C", "0.00", "93.46", "4.40", "0.0000", 10)
)

$r = 2
foreach ($row in $rows) {
    $q3.Cells.Item($r, 1).Value2 = $r - 2
    Set-TextValue $q3 $r 2 $row[0]
    Set-TextValue $q3 $r 3 $row[1]
    Set-TextValue $q3 $r 4 $row[2]
    Set-TextValue $q3 $r 5 $row[3]
    Set-TextValue $q3 $r 6 $row[4]
    Set-TextValue $q3 $r 7 $row[5]
    $q3.Cells.Item($r, 8).Value2 = $row[6]
    $r = $r + 1
}

# Last row's "持有市值(亿元)" (column G) is a literal number 0, not text.
$q3.Cells.Item(10, 7).Value2 = 0

# ---------------------------------------------------------------------
# Step 3: match the header/index-column styling used elsewhere in this
# workbook (bold + border + center/top alignment) by copying the format
# already applied to the "总计" sheet's own header cells.
# ---------------------------------------------------------------------
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$total.Range("A2").Copy()
$q3.Range("A2:A10").PasteSpecial(-4122)

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# Step 4: update the "总计" sheet - insert a new summary row for
# 2022-Q3 above the existing 2022-Q2 row, pushing it down to row 3.
# ---------------------------------------------------------------------
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)

$oldB2 = $total.Range("B2").Value2
$oldC2 = $total.Range("C2").Value2
$oldD2 = $total.Range("D2").Value2

$total.Range("A3").Value2 = 1
$total.Range("B3").Value2 = $oldB2
$total.Range("C3").Value2 = $oldC2
$total.Range("D3").Value2 = $oldD2

$total.Range("A2").Value2 = 0
$total.Range("B2").Value2 = "2022-Q3"
$total.Range("C2").Value2 = 9
$total.Range("D2").Value2 = 0.08

$total.Range("A1").Select()
